$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.043815935548006
$ws.Range("D2").Value = 1.049538552686424
$ws.Range("E2").Value = 1.051145469825808
$ws.Range("F2").Value = 1.060444702202214
$ws.Range("I2").Value = 1.042099601595074
$ws.Range("J2").Value = 1.048884356127076
$ws.Range("K2").Value = 1.052295331454685
$ws.Range("L2").Value = 1.053897783459321
$ws.Range("M2").Value = 1.063171463885072
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.044651734251337
$ws.Range("D3").Value = 1.050179079792162
$ws.Range("E3").Value = 1.051877682932876
$ws.Range("F3").Value = 1.061230918130278
$ws.Range("I3").Value = 1.042279994197488
$ws.Range("J3").Value = 1.049367719817893
$ws.Range("K3").Value = 1.052748568350864
$ws.Range("L3").Value = 1.054442791257497
$ws.Range("M3").Value = 1.06377217776078
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.045193311793747
$ws.Range("D4").Value = 1.050594204827494
$ws.Range("E4").Value = 1.05235252093768
$ws.Range("F4").Value = 1.061740776097747
$ws.Range("I4").Value = 1.042395899665549
$ws.Range("J4").Value = 1.049680567186651
$ws.Range("K4").Value = 1.053041814022179
$ws.Range("L4").Value = 1.054795823223241
$ws.Range("M4").Value = 1.064161341732753
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.045421171120801
$ws.Range("D5").Value = 1.050768879800691
$ws.Range("E5").Value = 1.052552391564493
$ws.Range("F5").Value = 1.061955387334944
$ws.Range("I5").Value = 1.042444429123199
$ws.Range("J5").Value = 1.049812105696021
$ws.Range("K5").Value = 1.05316508620039
$ws.Range("L5").Value = 1.05494432621963
$ws.Range("M5").Value = 1.064325055170559
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.045459440215225
$ws.Range("D6").Value = 1.050798217614563
$ws.Range("E6").Value = 1.05258596525072
$ws.Range("F6").Value = 1.061991437099861
$ws.Range("I6").Value = 1.042452565847062
$ws.Range("J6").Value = 1.049834192575352
$ws.Range("K6").Value = 1.053185783614536
$ws.Range("L6").Value = 1.054969265662977
$ws.Range("M6").Value = 1.064352549714412
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.045196355755686
$ws.Range("D7").Value = 1.050596538230818
$ws.Range("E7").Value = 1.052355190643909
$ws.Range("F7").Value = 1.061743642697277
$ws.Range("I7").Value = 1.042396548894332
$ws.Range("J7").Value = 1.049682324743708
$ws.Range("K7").Value = 1.053043461225413
$ws.Range("L7").Value = 1.054797807182284
$ws.Range("M7").Value = 1.064163528853916
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.044098239213874
$ws.Range("D8").Value = 1.049754883794086
$ws.Range("E8").Value = 1.051392706671527
$ws.Range("F8").Value = 1.060710173647076
$ws.Range("I8").Value = 1.042160735583044
$ws.Range("J8").Value = 1.049047694019337
$ws.Range("K8").Value = 1.052448510140689
$ws.Range("L8").Value = 1.054081892730562
$ws.Range("M8").Value = 1.063374381295929
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.042169116813205
$ws.Range("D9").Value = 1.048276927857721
$ws.Range("E9").Value = 1.049704793788198
$ws.Range("F9").Value = 1.058897766640684
$ws.Range("I9").Value = 1.041738950584611
$ws.Range("J9").Value = 1.047930057482039
$ws.Range("K9").Value = 1.051399970041467
$ws.Range("L9").Value = 1.05282330346788
$ws.Range("M9").Value = 1.061987414684262
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.040887108618459
$ws.Range("D10").Value = 1.047295196450146
$ws.Range("E10").Value = 1.048585086532085
$ws.Range("F10").Value = 1.057695466579441
$ws.Range("I10").Value = 1.041453601914086
$ws.Range("J10").Value = 1.047185495002342
$ws.Range("K10").Value = 1.050700918160561
$ws.Range("L10").Value = 1.05198631457484
$ws.Range("M10").Value = 1.061065299510577
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.040332973678517
$ws.Range("D11").Value = 1.046870967050995
$ws.Range("E11").Value = 1.048101586469411
$ws.Range("F11").Value = 1.057176299875375
$ws.Range("I11").Value = 1.041329065857246
$ws.Range("J11").Value = 1.046863233978085
$ws.Range("K11").Value = 1.050398232498861
$ws.Range("L11").Value = 1.051624399658693
$ws.Range("M11").Value = 1.060666634663245
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.040127292808015
$ws.Range("D12").Value = 1.046713521522233
$ws.Range("E12").Value = 1.047922196452124
$ws.Range("F12").Value = 1.056983676296424
$ws.Range("I12").Value = 1.041282661409161
$ws.Range("J12").Value = 1.046743554239387
$ws.Range("K12").Value = 1.050285804244548
$ws.Range("L12").Value = 1.0514900460581
$ws.Range("M12").Value = 1.060518647230812
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.040171405257563
$ws.Range("D13").Value = 1.0467472881111
$ws.Range("E13").Value = 1.047960666981426
$ws.Range("F13").Value = 1.057024984806473
$ws.Range("I13").Value = 1.041292621925092
$ws.Range("J13").Value = 1.046769224915462
$ws.Range("K13").Value = 1.050309920342894
$ws.Range("L13").Value = 1.051518861819554
$ws.Range("M13").Value = 1.060550386735954
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.040315968972892
$ws.Range("D14").Value = 1.046857949847288
$ws.Range("E14").Value = 1.048086753872908
$ws.Range("F14").Value = 1.057160373089741
$ws.Range("I14").Value = 1.041325233029409
$ws.Range("J14").Value = 1.046853340749403
$ws.Range("K14").Value = 1.050388939081495
$ws.Range("L14").Value = 1.05161329235373
$ws.Range("M14").Value = 1.060654400039322
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.040405059316462
$ws.Range("D15").Value = 1.046926149765428
$ws.Range("E15").Value = 1.048164467185743
$ws.Range("F15").Value = 1.057243819237105
$ws.Range("I15").Value = 1.041345306455145
$ws.Range("J15").Value = 1.046905170295059
$ws.Range("K15").Value = 1.050437625531358
$ws.Range("L15").Value = 1.051671484467269
$ws.Range("M15").Value = 1.060718498634638
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.040923905597174
$ws.Range("D16").Value = 1.047323369585478
$ws.Range("E16").Value = 1.048617203256315
$ws.Range("F16").Value = 1.057729952425211
$ws.Range("I16").Value = 1.04146184641662
$ws.Range("J16").Value = 1.047206885462979
$ws.Range("K16").Value = 1.050721006708998
$ws.Range("L16").Value = 1.052010344487811
$ws.Range("M16").Value = 1.061091770732876
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.041249628643266
$ws.Range("D17").Value = 1.047572768512335
$ws.Range("E17").Value = 1.048901552973092
$ws.Range("F17").Value = 1.058035277205382
$ws.Range("I17").Value = 1.041534687493847
$ws.Range("J17").Value = 1.047396181740938
$ws.Range("K17").Value = 1.050898767485666
$ws.Range("L17").Value = 1.052223039220847
$ws.Range("M17").Value = 1.061326080896414
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.041439712114576
$ws.Range("D18").Value = 1.047718322188808
$ws.Range("E18").Value = 1.04906753856152
$ws.Range("F18").Value = 1.05821350650567
$ws.Range("I18").Value = 1.041577080029232
$ws.Range("J18").Value = 1.047506608417228
$ws.Range("K18").Value = 1.051002453056457
$ws.Range("L18").Value = 1.05234714925331
$ws.Range("M18").Value = 1.061462809476644
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.041504541678387
$ws.Range("D19").Value = 1.047767966335509
$ws.Range("E19").Value = 1.049124157212142
$ws.Range("F19").Value = 1.05827430152194
$ws.Range("I19").Value = 1.041591518748027
$ws.Range("J19").Value = 1.047544263283518
$ws.Range("K19").Value = 1.051037807235285
$ws.Range("L19").Value = 1.052389475795401
$ws.Range("M19").Value = 1.061509440416201
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.041214671831583
$ws.Range("D20").Value = 1.047546001716821
$ws.Range("E20").Value = 1.048871031560175
$ws.Range("F20").Value = 1.058002504402697
$ws.Range("I20").Value = 1.041526882097812
$ws.Range("J20").Value = 1.047375870649781
$ws.Range("K20").Value = 1.050879695352265
$ws.Range("L20").Value = 1.052200214030666
$ws.Range("M20").Value = 1.061300935474175
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.04027339441976
$ws.Range("D21").Value = 1.04682535906039
$ws.Range("E21").Value = 1.048049618793948
$ws.Range("F21").Value = 1.05712049857838
$ws.Range("I21").Value = 1.041315633902107
$ws.Range("J21").Value = 1.046828570104467
$ws.Range("K21").Value = 1.050365669947207
$ws.Range("L21").Value = 1.051585482757133
$ws.Range("M21").Value = 1.060623768092989
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.03968244192837
$ws.Range("D22").Value = 1.046373028006811
$ws.Range("E22").Value = 1.047534341994152
$ws.Range("F22").Value = 1.05656720914017
$ws.Range("I22").Value = 1.041181968215301
$ws.Range("J22").Value = 1.046484590560359
$ws.Range("K22").Value = 1.050042497964133
$ws.Range("L22").Value = 1.051199427211426
$ws.Range("M22").Value = 1.060198553353738
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.039995634330643
$ws.Range("D23").Value = 1.046612744042384
$ws.Range("E23").Value = 1.047807387636667
$ws.Range("F23").Value = 1.056860397944913
$ws.Range("I23").Value = 1.041252906822862
$ws.Range("J23").Value = 1.046666927830465
$ws.Range("K23").Value = 1.050213815540461
$ws.Range("L23").Value = 1.051404039342438
$ws.Range("M23").Value = 1.060423915280349
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.041230467017617
$ws.Range("D24").Value = 1.047558096221802
$ws.Range("E24").Value = 1.048884822473929
$ws.Range("F24").Value = 1.058017312593817
$ws.Range("I24").Value = 1.041530409312178
$ws.Range("J24").Value = 1.047385048316496
$ws.Range("K24").Value = 1.050888313226366
$ws.Range("L24").Value = 1.052210527601054
$ws.Range("M24").Value = 1.061312297424724
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.042667131082196
$ws.Range("D25").Value = 1.048658393413763
$ws.Range("E25").Value = 1.050140186920195
$ws.Range("F25").Value = 1.059365273859582
$ws.Range("I25").Value = 1.041848728390913
$ws.Range("J25").Value = 1.048218906175673
$ws.Range("K25").Value = 1.051671053271619
$ws.Range("L25").Value = 1.053148320704819
$ws.Range("M25").Value = 1.062345540366934

Write-Output "Applied 240 cell updates"